$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.510597467422485
$ws.Range("B1").Value = 2.02187967300415
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.45405375957489
$ws.Range("E1").Value = 0.6595985293388367
